$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Doctor's Details")

$ws.Range("A2").Value = "Dr. Dinesh Shetty"
$ws.Range("B2").Value = "BDS`nDentist, Dental Surgeon`n20 Years Experience Overall"

$ws.Range("A3").Value = "Dr. Lokesh Babu"
$ws.Range("B3").Value = "BDS, MDS - Oral & Maxillofacial Surgery`nDentist, Oral And MaxilloFacial Surgeon, Implantologist`n25 Years Experience Overall  (24 years as specialist)"

$ws.Range("A4").Value = "Dr. Narayan Babu"
$ws.Range("B4").Value = "BDS`nDentist, Dental Surgeon, Cosmetic/Aesthetic Dentist`n16 Years Experience Overall"

$ws.Range("A5").Value = "Dr. Divya Puranik"
$ws.Range("B5").Value = "BDS, MDS - Orthodontics and Dentofacial Orthopaedics`nOrthodontist, Dentofacial Orthopedist`n15 Years Experience Overall  (9 years as specialist)"

$ws.Range("A6").Value = "Dr. Nikhar Ravinder"
$ws.Range("B6").Value = "BDS, MDS - Orthodontics and Dentofacial Orthopaedics`nOrthodontist, Dentist, Dental Surgeon, Dentofacial Orthopedist`n22 Years Experience Overall  (17 years as specialist)"
